$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 7, shifting existing rows 7-50 down to 8-51.
$ws.Rows.Item(7).Insert(-4121)

# Fill in the new row 7 with data for a new weekly price observation.
$ws.Cells.Item(7, 1).Value = 2
$ws.Cells.Item(7, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(7, 3).Value = "Coquimbo"
$ws.Cells.Item(7, 4).Value = 44532
$ws.Cells.Item(7, 5).Value = 4
$ws.Cells.Item(7, 6).Value = 100112030
$ws.Cells.Item(7, 7).Value = "Poroto granado"
$ws.Cells.Item(7, 8).Value = "Sin especificar"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 400
$ws.Cells.Item(7, 11).Value = 28000
$ws.Cells.Item(7, 12).Value = 30000
$ws.Cells.Item(7, 13).Value = 29000
$ws.Cells.Item(7, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(7, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(7, 16).Value = 1933
$ws.Cells.Item(7, 17).Value = 15
$ws.Cells.Item(7, 18).Value = "Hortaliza"
